# Weekly refresh of the Rabanito (Vega Modelo de Temuco) price series:
# three new weekly records are inserted into the data block (which runs
# from row 13 to row 25), pushing the existing rows down and extending
# the sheet to row 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common/constant columns shared by every row in this sub-sheet.
$colA = 10
$colB = "Vega Modelo de Temuco"
$colC = "La Araucanía"
$colE = 9
$colF = 300000001
$colG = "Rabanito"
$colH = "Sin especificar"
$colI = "Primera"
$colN = "`$/docena de paquetes"
$colQ = 12
$colR = "Hortaliza"

function Set-DataRow {
    param($rowNum, $dateSerial, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg)

    $ws.Cells.Item($rowNum, 1).Value2 = $colA
    $ws.Cells.Item($rowNum, 2).Value2 = $colB
    $ws.Cells.Item($rowNum, 3).Value2 = $colC
    $ws.Cells.Item($rowNum, 4).Value2 = $dateSerial
    $ws.Cells.Item($rowNum, 5).Value2 = $colE
    $ws.Cells.Item($rowNum, 6).Value2 = $colF
    $ws.Cells.Item($rowNum, 7).Value2 = $colG
    $ws.Cells.Item($rowNum, 8).Value2 = $colH
    $ws.Cells.Item($rowNum, 9).Value2 = $colI
    $ws.Cells.Item($rowNum, 10).Value2 = $volumen
    $ws.Cells.Item($rowNum, 11).Value2 = $precioMin
    $ws.Cells.Item($rowNum, 12).Value2 = $precioMax
    $ws.Cells.Item($rowNum, 13).Value2 = $precioProm
    $ws.Cells.Item($rowNum, 14).Value2 = $colN
    $ws.Cells.Item($rowNum, 15).Value2 = $origen
    $ws.Cells.Item($rowNum, 16).Value2 = $precioKg
    $ws.Cells.Item($rowNum, 17).Value2 = $colQ
    $ws.Cells.Item($rowNum, 18).Value2 = $colR
}

# 1) New record (2021-08-20) becomes the new first data row, row 13;
#    everything that was row 13-25 shifts down to 14-26.
$ws.Range("A13").EntireRow.Insert()
Set-DataRow 13 44428 10 7000 7000 7000 "Provincia de Cautín" 583

# 2) New record (2021-08-27) is inserted right after the row that now
#    holds the old row-16 data (2021-08-02), i.e. at row 18; rows
#    18-26 shift down to 19-27.
$ws.Range("A18").EntireRow.Insert()
Set-DataRow 18 44435 30 7000 7000 7000 "Provincia de Cautín" 583

# 3) New record (2021-08-24) is appended at the very end, row 28.
$ws.Range("A28").EntireRow.Insert()
Set-DataRow 28 44432 30 7000 7000 7000 "Provincia de Cautín" 583
